$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 202, pushing existing rows 202..324 down to 203..325
$ws.Rows.Item(202).Insert()

# Populate the newly inserted row 202 with the new data record.
# Columns A, B, C, E, F, G, H, I, J, N, Q, R are identical to the row that
# used to occupy 202 (now at 203), so they're already correct after the
# insert copied formatting/values down - only set the cells that actually
# carry new data for this record, plus the unchanged-but-needed ones that
# the Insert() left blank.
$ws.Range("A202").Value = 3
$ws.Range("B202").Value = "Femacal de La Calera"
$ws.Range("C202").Value = "Coquimbo"
$ws.Range("D202").Value = 44767
$ws.Range("E202").Value = 5
$ws.Range("F202").Value = 100112001
$ws.Range("G202").Value = "Berenjena"
$ws.Range("H202").Value = "Sin especificar"
$ws.Range("I202").Value = "Primera"
$ws.Range("J202").Value = 80
$ws.Range("K202").Value = 9000
$ws.Range("L202").Value = 9000
$ws.Range("M202").Value = 9000
$ws.Range("N202").Value = "$/caja 60 unidades"
$ws.Range("O202").Value = "Provincia de Limarí"
$ws.Range("P202").Value = 150
$ws.Range("Q202").Value = 60
$ws.Range("R202").Value = "Hortaliza"
